# Distribution list.docx edit
#  - "Added support for different styles of pagination":
#      flip the section's page setup from portrait (Letter, 8.5x11)
#      to landscape (11x8.5) and swap/adjust the margins to match.
#  - "Fixed command line operation to use all fields":
#      the stray _GoBack bookmark (left over from the last cursor
#      position when the file was last saved) now sits on the very
#      first paragraph instead of next to "Notes:" further down.

$d = $word.ActiveDocument

# --- Move the _GoBack bookmark to the first paragraph of the document ---
# Word only ever keeps a single "_GoBack" bookmark; adding it anew moves
# it (removing the previous bookmarkStart/bookmarkEnd pair wherever it was).
$firstPara = $d.Paragraphs(1)
$d.Bookmarks.Add("_GoBack", $firstPara.Range)

# --- Switch the page from portrait to landscape and rebalance margins ---
$ps = $d.PageSetup
$ps.Orientation = 1            # wdOrientLandscape
$ps.PageWidth = 792            # 15840 twips
$ps.PageHeight = 612            # 12240 twips
$ps.TopMargin = 72             # 1440 twips
$ps.RightMargin = 43.2         # 864 twips
$ps.BottomMargin = 72          # 1440 twips
$ps.LeftMargin = 36            # 720 twips
